# Sushrut/currentSchedule.xlsx - "added routes 4,5,6 and made f function in main"
#
# Logical summary of the edit:
#   1. A new (blank) column is inserted before column H, pushing the old
#      "StartTimingsfromthestartingpoint" block (H:U) one column to the right
#      (I:V), and re-titling that header with normal spacing.
#   2. Two new helper columns are added further right:
#        X -> "Frequency from 8:10 to 8:45" (trip counts per route/bus)
#        Z -> "Distance Travelled in meters" (=count * per-trip distance)
#   3. Row 18 sums up the new distance column.
#   4. The named range / used dimension grows by one column (U->V) and two
#      rows (16->18) to cover the new data.
#   5. Minor cosmetic bits: page setup, active selection cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a blank column before H (shifts H:CQ -> I:CR) -----------
$ws.Columns("H").Insert()

# Re-title the (now shifted) header cell with normal word spacing.
$ws.Range("I1").Value = "Start Timings from the starting point"

# --- 2. New helper columns: trip-frequency (X) and distance (Z) --------
$ws.Range("X1").Value = "Frequency from 8:10 to 8:45"
$ws.Range("Z1").Value = "Distance Travelled in meters"

# Route frequency counts (routes 1-3 use the 1300m leg, 4-5 the 1170m leg,
# 6-8 the 990m+1170m legs, 9-11 the 880m leg).
$ws.Range("X2").Value = 4
$ws.Range("X3:X5").Value = 4
$ws.Range("X5").Value = 3
$ws.Range("X6").Value = 4
$ws.Range("X7").Value = 4
$ws.Range("X8").Value = 4
$ws.Range("X9:X13").Value = 4
$ws.Range("X14").Value = 4
$ws.Range("X15:X16").Value = 4

# General (non time-formatted) display for the newly-typed route counts.
$ws.Range("X6:X16").NumberFormat = "General"

# Distance travelled = frequency * round-trip distance for that leg.
$ws.Range("Z2").Formula = "=X2*1300*2"
$ws.Range("Z3:Z5").Formula = "=X3*1300*2"
$ws.Range("Z6").Formula = "=X6*1170*2"
$ws.Range("Z7").Formula = "=X7*1170*2"
$ws.Range("Z8").Formula = "=X8*(990 + 1170)"
$ws.Range("Z9:Z13").Formula = "=X9*(990 + 1170)"
$ws.Range("Z14").Formula = "=X14*880*2"
$ws.Range("Z15:Z16").Formula = "=X15*880*2"

# --- 3. Grand total of the new distance column -------------------------
$ws.Range("Z18").Formula = "=SUM(Z2:Z17)"

# --- 4. Keep the named range / dimension in sync with the extra column -
$wb.Names.Item("tumtum_sche_1").RefersTo = "=Sheet1!`$A`$1:`$V`$16"

# --- 5. Cosmetics: page setup + active selection ------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
[void]$ws.Range("Z18").Select()
